$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.294.67"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "3.829.72"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.01"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.60"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").Value = "3.827.71"
$ws.Range("E7").Value = "  -1.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("E10").Value = "  -2.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.47"
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("E12").Value = "  -2.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000268"
$ws.Range("E13").Value = "  +4.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.09"
$ws.Range("E14").Value = "  -3.34%  "
$ws.Range("D15").Value = "4.473.39"
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").Value = "3.826.39"
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("D17").Value = "68.328.40"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.49"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.38"
$ws.Range("E19").Value = "  -3.83%  "
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.06"
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.33"
$ws.Range("E22").Value = "  -4.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.734"
$ws.Range("E23").Value = "  -2.20%  "
$ws.Range("E24").Value = "  -3.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.12"
$ws.Range("E25").Value = "  -2.67%  "
$ws.Range("E26").Value = "  -3.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.17"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").Value = "3.980.22"
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.68"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.57"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("E34").Value = "  -4.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.38"
$ws.Range("E35").Value = "  -2.62%  "
$ws.Range("D36").Value = "3.795.33"
$ws.Range("E36").Value = "  -1.58%  "
$ws.Range("E37").Value = "  -2.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.68"
$ws.Range("E38").Value = "  +10.67%  "
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.94"
$ws.Range("E41").Value = "  -3.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.314"
$ws.Range("E43").Value = "  -4.57%  "
$ws.Range("E44").Value = "  -5.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.73"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "416.22"
$ws.Range("E47").Value = "  -5.26%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.16"
$ws.Range("E48").Value = "  -2.28%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000291"
$ws.Range("E49").Value = "  +5.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.89"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("E51").Value = "  -2.70%  "
